$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 415.58334
$ws.Range("I33").Value = 407.9091
$ws.Range("K33").Value = 407.9091
$ws.Range("M33").Value = -178.9091
# Row 80
$ws.Range("H80").Value = 847.5
$ws.Range("J80").Value = 847.5
$ws.Range("L80").Value = 2542.5
$ws.Range("N80").Value = -4538.5
# Row 83
$ws.Range("H83").Value = 847.5
$ws.Range("J83").Value = 847.5
$ws.Range("L83").Value = 7627.5
$ws.Range("N83").Value = -17611.5
# Row 115
$ws.Range("H115").Value = 10449.75
$ws.Range("I115").Value = 8933
$ws.Range("J115").Value = 15000
$ws.Range("K115").Value = 26799
$ws.Range("L115").Value = 45000
$ws.Range("M115").Value = -25232
$ws.Range("N115").Value = -48134
# Row 118
$ws.Range("H118").Value = 2963.5908
$ws.Range("J118").Value = 2979.95
$ws.Range("L118").Value = 8939.849999999999
$ws.Range("N118").Value = -12253.85
# Row 125
$ws.Range("H125").Value = 1805.3334
$ws.Range("I125").Value = 1199.5
$ws.Range("K125").Value = 10795.5
$ws.Range("M125").Value = -8335.5
# Row 127
$ws.Range("H127").Value = 1993.4
$ws.Range("I127").Value = 1992.6666
$ws.Range("J127").Value = 1994.5
$ws.Range("K127").Value = 5977.9998
$ws.Range("L127").Value = 5983.5
$ws.Range("M127").Value = -1017.9998
$ws.Range("N127").Value = -15903.5
# Row 129
$ws.Range("H129").Value = 2098.7144
$ws.Range("J129").Value = 1500
$ws.Range("L129").Value = 4500
$ws.Range("N129").Value = -14500
# Row 135
$ws.Range("H135").Value = 1378.4286
$ws.Range("I135").Value = 949.8
$ws.Range("K135").Value = 8548.199999999999
$ws.Range("M135").Value = -6013.199999999999
# Row 137
$ws.Range("H137").Value = 2422.8857
$ws.Range("I137").Value = 1458.45
$ws.Range("K137").Value = 4375.35
$ws.Range("M137").Value = -1825.35
# Row 138
$ws.Range("H138").Value = 4296.148
$ws.Range("J138").Value = 4534.067
$ws.Range("L138").Value = 13602.201
$ws.Range("N138").Value = -23882.201

$ws = $wb.Worksheets.Item("ARM")
# Row 7
$ws.Range("H7").Value = 2000
$ws.Range("I7").Value = 2000
$ws.Range("K7").Value = 2000
$ws.Range("M7").Value = -1888
# Row 32
$ws.Range("H32").Value = 6459.675
$ws.Range("I32").Value = 6201.9624
$ws.Range("J32").Value = 13332
$ws.Range("K32").Value = 6201.9624
$ws.Range("L32").Value = 13332
$ws.Range("M32").Value = -5914.9624
$ws.Range("N32").Value = -13906
# Row 126
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
# Row 132
$ws.Range("H132").Value = 2823.3845
$ws.Range("I132").Value = 2024.2941
$ws.Range("J132").Value = 4332.778
$ws.Range("K132").Value = 6072.8823
$ws.Range("L132").Value = 12998.334
$ws.Range("M132").Value = -3542.8823
$ws.Range("N132").Value = -18058.334

$ws = $wb.Worksheets.Item("BSM")
# Row 122
$ws.Range("H122").Value = 10000
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
# Row 134
$ws.Range("H134").Value = 2285
$ws.Range("I134").Value = 1888.9286
$ws.Range("J134").Value = 3869.2856
$ws.Range("K134").Value = 5666.7858
$ws.Range("L134").Value = 11607.8568
$ws.Range("M134").Value = -3131.7858
$ws.Range("N134").Value = -16677.8568

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 10004
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
# Row 31
$ws.Range("H31").Value = 2507.8572
$ws.Range("I31").Value = 2563.3333
$ws.Range("J31").Value = 2408
$ws.Range("K31").Value = 2563.3333
$ws.Range("L31").Value = 2408
$ws.Range("M31").Value = -2268.3333
$ws.Range("N31").Value = -2998
# Row 34
$ws.Range("H34").Value = 2507.8572
$ws.Range("I34").Value = 2563.3333
$ws.Range("J34").Value = 2408
$ws.Range("K34").Value = 2563.3333
$ws.Range("L34").Value = 2408
$ws.Range("M34").Value = -2361.3333
$ws.Range("N34").Value = -2812
# Row 41
$ws.Range("H41").Value = 38991.25
$ws.Range("J41").Value = 38991.25
$ws.Range("L41").Value = 38991.25
$ws.Range("N41").Value = -39847.25
# Row 42
$ws.Range("H42").Value = 2056
$ws.Range("I42").Value = 2056
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 2056
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -1463
$ws.Range("N42").ClearContents()
# Row 141
$ws.Range("H141").Value = 50448.855
$ws.Range("J141").Value = 50448.855
$ws.Range("L141").Value = 50448.855
$ws.Range("N141").Value = -60808.855

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 4
$ws.Range("H4").Value = 22487870
$ws.Range("I4").Value = 33232988
$ws.Range("K4").Value = 99698964
$ws.Range("M4").Value = -99698852
# Row 137
$ws.Range("H137").Value = 4992.636
$ws.Range("I137").Value = 3950
$ws.Range("J137").Value = 5224.3335
$ws.Range("K137").Value = 11850
$ws.Range("L137").Value = 15673.0005
$ws.Range("M137").Value = -6750
$ws.Range("N137").Value = -25873.0005
# Row 140
$ws.Range("H140").Value = 1162.8823
$ws.Range("I140").Value = 1162.8823
$ws.Range("K140").Value = 3488.6469
$ws.Range("M140").Value = 1691.3531

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 553.9
$ws.Range("I107").Value = 556.25
$ws.Range("J107").Value = 544.5
$ws.Range("K107").Value = 556.25
$ws.Range("L107").Value = 544.5
$ws.Range("M107").Value = 1363.75
$ws.Range("N107").Value = -4384.5
# Row 132
$ws.Range("H132").Value = 4990.4116
$ws.Range("I132").Value = 5565.5557
$ws.Range("K132").Value = 16696.6671
$ws.Range("M132").Value = -14166.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4096
$ws.Range("I61").Value = 4248.6
$ws.Range("J61").Value = 3841.6667
$ws.Range("K61").Value = 4248.6
$ws.Range("L61").Value = 3841.6667
$ws.Range("M61").Value = -4046.6
$ws.Range("N61").Value = -4245.6667
# Row 113
$ws.Range("H113").Value = 4096
$ws.Range("I113").Value = 4248.6
$ws.Range("J113").Value = 3841.6667
$ws.Range("K113").Value = 4248.6
$ws.Range("L113").Value = 3841.6667
$ws.Range("M113").Value = -2078.6
$ws.Range("N113").Value = -8181.6667
# Row 132
$ws.Range("H132").Value = 5266.3335
$ws.Range("I132").Value = 5057
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 15171
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -12641
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 71333
$ws.Range("J46").Value = 71999.5
$ws.Range("L46").Value = 71999.5
$ws.Range("N46").Value = -72461.5
# Row 107
$ws.Range("H107").Value = 561.8333
$ws.Range("I107").Value = 514.2
$ws.Range("K107").Value = 1542.6
$ws.Range("M107").Value = 377.3999999999999
# Row 134
$ws.Range("H134").Value = 71333
$ws.Range("J134").Value = 71999.5
$ws.Range("L134").Value = 215998.5
$ws.Range("N134").Value = -221068.5
